# Deploying to gh-pages from @ Alvearie/alvearie-fhir-ig@8e4a450c507ef6f746e072652acbb72e9504f19a
# Update the FHIR StructureDefinition export workbook:
#  - bump Version/Date on the "Metadata" sheet
#  - replace the duplicated "Contact / No display for ContactDetail" rows
#    with real "Publisher"/"Jurisdiction" values, removing the duplicate row
#  - give the root Extension row on "Elements" its real Short/Definition text

$wb = $excel.ActiveWorkbook

$meta = $wb.Worksheets.Item("Metadata")

# Version: 5.0.0 -> 6.0.0
$meta.Range("B3").Value = "6.0.0"

# Date: refresh the IG build timestamp
$meta.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher value was blank ("Publisher" label with nothing next to it) - fill it in
$meta.Range("B9").Value = "Alvearie Team"

# The old row 10 ("Contact" / "No display for ContactDetail") becomes the
# Jurisdiction row
$meta.Range("A10").Value = "Jurisdiction"
$meta.Range("B10").Value = "United States of America"

# Row 11 used to be an exact duplicate of the old row 10 ("Contact" /
# "No display for ContactDetail"); delete it so everything below shifts up
# one row (dimension goes from A1:B21 to A1:B20)
$meta.Rows.Item(11).Delete()

$elements = $wb.Worksheets.Item("Elements")

# Root Extension row (row 2): give it the profile-specific Short/Definition
# text instead of the generic "Extension" / "An Extension"
$elements.Range("K2").Value = "Shortterm Care Coverage Indicator"
$elements.Range("L2").Value = "Indicates whether the member or employee has short-term care benefit coverage"
